$d = $word.ActiveDocument

# 1. Replace the title text
$d.Content.Find.Execute("2.2 - Debate I", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Placeholder - Check Back Later", 2)

# 2. Remove the trailing " :::" runs after "...general edification later."
$d.Content.Find.Execute(" :::", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
